$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1641.7059
$ws.Range("I6").Value = 131.125
$ws.Range("K6").Value = 393.375
$ws.Range("M6").Value = -281.375
$ws.Range("H33").Value = 824.85
$ws.Range("I33").Value = 965.6667
$ws.Range("J33").Value = 709.63635
$ws.Range("K33").Value = 965.6667
$ws.Range("L33").Value = 709.63635
$ws.Range("M33").Value = -736.6667
$ws.Range("N33").Value = -1167.63635
$ws.Range("H58").Value = 51405.9
$ws.Range("I58").Value = 381.625
$ws.Range("J58").Value = 85422.086
$ws.Range("K58").Value = 1144.875
$ws.Range("L58").Value = 256266.258
$ws.Range("M58").Value = -994.875
$ws.Range("N58").Value = -256566.258
$ws.Range("H98").Value = 2267.9333
$ws.Range("I98").Value = 2233.5454
$ws.Range("J98").Value = 2362.5
$ws.Range("K98").Value = 2233.5454
$ws.Range("L98").Value = 2362.5
$ws.Range("M98").Value = -735.5454
$ws.Range("N98").Value = -5358.5
$ws.Range("H122").Value = 2267.9333
$ws.Range("I122").Value = 2233.5454
$ws.Range("J122").Value = 2362.5
$ws.Range("K122").Value = 6700.6362
$ws.Range("L122").Value = 7087.5
$ws.Range("M122").Value = -4250.6362
$ws.Range("N122").Value = -11987.5
$ws.Range("H138").Value = 1939.6111
$ws.Range("I138").Value = 1335.9744
$ws.Range("J138").Value = 2401.2156
$ws.Range("K138").Value = 4007.9232
$ws.Range("L138").Value = 7203.6468
$ws.Range("M138").Value = 1132.0768
$ws.Range("N138").Value = -17483.6468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8104.3765
$ws.Range("I32").Value = 6687.4863
$ws.Range("J32").Value = 17636.182
$ws.Range("K32").Value = 6687.4863
$ws.Range("L32").Value = 17636.182
$ws.Range("M32").Value = -6400.4863
$ws.Range("N32").Value = -18210.182
$ws.Range("H45").Value = 1466.7435
$ws.Range("I45").Value = 1078.6875
$ws.Range("J45").Value = 3240.7144
$ws.Range("K45").Value = 1078.6875
$ws.Range("L45").Value = 3240.7144
$ws.Range("M45").Value = -701.6875
$ws.Range("N45").Value = -3994.7144
$ws.Range("H80").Value = 25203.166
$ws.Range("J80").Value = 25203.166
$ws.Range("L80").Value = 25203.166
$ws.Range("N80").Value = -27199.166
$ws.Range("H83").Value = 25203.166
$ws.Range("J83").Value = 25203.166
$ws.Range("L83").Value = 75609.49800000001
$ws.Range("N83").Value = -85593.49800000001
$ws.Range("H122").Value = 2525.739
$ws.Range("I122").Value = 1830.8125
$ws.Range("J122").Value = 4114.143
$ws.Range("K122").Value = 5492.4375
$ws.Range("L122").Value = 12342.429
$ws.Range("M122").Value = -3042.4375
$ws.Range("N122").Value = -17242.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4328.3794
$ws.Range("I31").Value = 2827.8948
$ws.Range("J31").Value = 7179.3
$ws.Range("K31").Value = 2827.8948
$ws.Range("L31").Value = 7179.3
$ws.Range("M31").Value = -2532.8948
$ws.Range("N31").Value = -7769.3
$ws.Range("H34").Value = 4328.3794
$ws.Range("I34").Value = 2827.8948
$ws.Range("J34").Value = 7179.3
$ws.Range("K34").Value = 2827.8948
$ws.Range("L34").Value = 7179.3
$ws.Range("M34").Value = -2625.8948
$ws.Range("N34").Value = -7583.3
$ws.Range("H107").Value = 1536.3
$ws.Range("I107").Value = 606.25
$ws.Range("J107").Value = 5256.5
$ws.Range("K107").Value = 606.25
$ws.Range("L107").Value = 5256.5
$ws.Range("M107").Value = 1313.75
$ws.Range("N107").Value = -9096.5
$ws.Range("H132").Value = 2261.8386
$ws.Range("I132").Value = 1635.5217
$ws.Range("J132").Value = 4062.5
$ws.Range("K132").Value = 4906.5651
$ws.Range("L132").Value = 12187.5
$ws.Range("M132").Value = -2376.5651
$ws.Range("N132").Value = -17247.5
$ws.Range("H134").Value = 13997.031
$ws.Range("I134").Value = 16620.2
$ws.Range("J134").Value = 4628.5713
$ws.Range("K134").Value = 49860.60000000001
$ws.Range("L134").Value = 13885.7139
$ws.Range("M134").Value = -47325.60000000001
$ws.Range("N134").Value = -18955.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3601.7144
$ws.Range("I70").Value = 2606
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 7818
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -7503
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 3601.7144
$ws.Range("I73").Value = 2606
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 7818
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -6726
$ws.Range("N73").Value = -14184
$ws.Range("H114").Value = 703.8461
$ws.Range("I114").Value = 215
$ws.Range("J114").Value = 2333.3333
$ws.Range("K114").Value = 645
$ws.Range("L114").Value = 6999.999899999999
$ws.Range("M114").Value = 2609
$ws.Range("N114").Value = -13507.9999
$ws.Range("H122").Value = 1344.8334
$ws.Range("I122").Value = 590
$ws.Range("J122").Value = 2099.6667
$ws.Range("K122").Value = 5310
$ws.Range("L122").Value = 18897.0003
$ws.Range("M122").Value = -2860
$ws.Range("N122").Value = -23797.0003
$ws.Range("H129").Value = 41717.69
$ws.Range("I129").Value = 6257.5
$ws.Range("J129").Value = 57477.777
$ws.Range("K129").Value = 18772.5
$ws.Range("L129").Value = 172433.331
$ws.Range("M129").Value = -13772.5
$ws.Range("N129").Value = -182433.331
$ws.Range("H131").Value = 1328.826
$ws.Range("J131").Value = 1105.7213
$ws.Range("L131").Value = 3317.1639
$ws.Range("N131").Value = -13397.1639

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 22172
$ws.Range("J95").Value = 22172
$ws.Range("L95").Value = 22172
$ws.Range("N95").Value = -27664
$ws.Range("H122").Value = 4921.3335
$ws.Range("I122").Value = 2635.3333
$ws.Range("J122").Value = 6445.3335
$ws.Range("K122").Value = 7905.999899999999
$ws.Range("L122").Value = 19336.0005
$ws.Range("M122").Value = -5455.999899999999
$ws.Range("N122").Value = -24236.0005
$ws.Range("H132").Value = 2494.742
$ws.Range("I132").Value = 2050.628
$ws.Range("J132").Value = 3499.842
$ws.Range("K132").Value = 6151.884
$ws.Range("L132").Value = 10499.526
$ws.Range("M132").Value = -3621.884
$ws.Range("N132").Value = -15559.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 7013.6665
$ws.Range("I45").Value = 3041
$ws.Range("K45").Value = 3041
$ws.Range("M45").Value = -2634
$ws.Range("H132").Value = 3682.15
$ws.Range("I132").Value = 2484.3
$ws.Range("J132").Value = 4880
$ws.Range("K132").Value = 7452.900000000001
$ws.Range("L132").Value = 14640
$ws.Range("M132").Value = -4922.900000000001
$ws.Range("N132").Value = -19700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1163.6786
$ws.Range("I113").Value = 381.14285
$ws.Range("J113").Value = 3511.2856
$ws.Range("K113").Value = 1143.42855
$ws.Range("L113").Value = 10533.8568
$ws.Range("M113").Value = 1026.57145
$ws.Range("N113").Value = -14873.8568
$ws.Range("H132").Value = 13495.866
$ws.Range("I132").Value = 2099.7273
$ws.Range("J132").Value = 44835.25
$ws.Range("K132").Value = 6299.1819
$ws.Range("L132").Value = 134505.75
$ws.Range("M132").Value = -3769.1819
$ws.Range("N132").Value = -139565.75
